$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "49.440.50"
$ws.Range("E2").Value = "  -0.93%  "
$ws.Range("D3").Value = "2.631.30"
$ws.Range("E3").Value = "  -0.81%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "111.78"
$ws.Range("E5").Value = "  +0.65%  "
$ws.Range("D6").Value = "324.49"
$ws.Range("E6").Value = "  -0.99%  "
$ws.Range("E7").Value = "  -1.27%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  -2.79%  "
$ws.Range("D10").Value = "39.55"
$ws.Range("E10").Value = "  -3.38%  "
$ws.Range("D11").Value = "19.82"
$ws.Range("E11").Value = "  -3.26%  "
$ws.Range("D12").Value = "0.0809"
$ws.Range("E12").Value = "  -1.50%  "
$ws.Range("E13").Value = "  +1.44%  "
$ws.Range("D14").Value = "7.32"
$ws.Range("E14").Value = "  -0.23%  "
$ws.Range("D15").Value = "3.041.30"
$ws.Range("E15").Value = "  -1.07%  "
$ws.Range("D16").Value = "2.621.34"
$ws.Range("E16").Value = "  -2.47%  "
$ws.Range("E17").Value = "  -3.90%  "
$ws.Range("D18").Value = "49.342.08"
$ws.Range("D19").Value = "12.86"
$ws.Range("E19").Value = "  -3.45%  "
$ws.Range("D20").Value = "2.93"
$ws.Range("E20").Value = "  -0.77%  "
$ws.Range("E21").Value = "  -2.20%  "
$ws.Range("D22").Value = "0.0₃0944"
$ws.Range("E22").Value = "  -1.74%  "
$ws.Range("D23").Value = "269.68"
$ws.Range("E23").Value = "  -3.95%  "
$ws.Range("D24").Value = "68.85"
$ws.Range("E24").Value = "  -5.86%  "
$ws.Range("E25").Value = "  -2.61%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.20"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.22%  "
$ws.Range("E27").Value = "  +0.05%  "
$ws.Range("D28").Value = "10.26"
$ws.Range("E28").Value = "  +2.82%  "
$ws.Range("E29").Value = "  -1.27%  "
$ws.Range("D30").Value = "0.136"
$ws.Range("E30").Value = "  -5.07%  "
$ws.Range("D31").Value = "34.56"
$ws.Range("E31").Value = "  -6.26%  "
$ws.Range("D32").Value = "49.45"
$ws.Range("E32").Value = "  -0.63%  "
$ws.Range("D33").Value = "5.49"
$ws.Range("E33").Value = "  +0.93%  "
$ws.Range("D34").Value = "0.0814"
$ws.Range("E34").Value = "  +1.58%  "
$ws.Range("E35").Value = "  -0.26%  "
$ws.Range("D36").Value = "18.87"
$ws.Range("E36").Value = "  -3.69%  "
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D37").Value = "4.88"
$ws.Range("E37").Value = "  +2.24%  "
$ws.Range("B38").Value = "ARBITRUM"
$ws.Range("C38").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D38").Value = "2.04"
$ws.Range("E38").Value = "  -1.68%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.10"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.98%  "
$ws.Range("D40").Value = "128.42"
$ws.Range("E40").Value = "  +1.72%  "
$ws.Range("E41").Value = "  -1.63%  "
$ws.Range("D42").Value = "22.14"
$ws.Range("E42").Value = "  -1.12%  "
$ws.Range("D43").Value = "0.0325"
$ws.Range("E43").Value = "  +3.47%  "
$ws.Range("D44").Value = "2.14"
$ws.Range("E44").Value = "  -4.37%  "
$ws.Range("D45").Value = "2.055.01"
$ws.Range("E45").Value = "  -0.86%  "
$ws.Range("E46").Value = "  -5.42%  "
$ws.Range("E47").Value = "  +6.30%  "
$ws.Range("E48").Value = "  -5.61%  "
$ws.Range("D49").Value = "8.89"
$ws.Range("E49").Value = "  -2.12%  "
$ws.Range("D50").Value = "59.01"
$ws.Range("E50").Value = "  +1.83%  "
$ws.Range("E51").Value = "  -4.26%  "
